# Refresh the "想去人数" (interest-count) column F across all four sheets
# to match a newer scrape of show.bilibili.com (gh-pages data regeneration).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 14481
$ws.Range("F4").Value = 14666
$ws.Range("F5").Value = 1378
$ws.Range("F7").Value = 5972
$ws.Range("F8").Value = 1000
$ws.Range("F15").Value = 2137
$ws.Range("F16").Value = 1246
$ws.Range("F18").Value = 922
$ws.Range("F19").Value = 39
$ws.Range("F20").Value = 2310
$ws.Range("F21").Value = 584
$ws.Range("F22").Value = 841
$ws.Range("F23").Value = 3427
$ws.Range("F26").Value = 2497
$ws.Range("F27").Value = 621
$ws.Range("F31").Value = 1097
$ws.Range("F32").Value = 1473
$ws.Range("F35").Value = 5057
$ws.Range("F36").Value = 4982
$ws.Range("F39").Value = 694
$ws.Range("F40").Value = 704
$ws.Range("F41").Value = 3332
$ws.Range("F45").Value = 130
$ws.Range("F47").Value = 4460
$ws.Range("F48").Value = 641
$ws.Range("F49").Value = 309

# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F20").Value = 19
$ws.Range("F27").Value = 1

# Sheet 3: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 7757
$ws.Range("F4").Value = 944

# Sheet 4: 全部类型 (All Types - combined view)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 7757
$ws.Range("F5").Value = 944
$ws.Range("F7").Value = 14481
$ws.Range("F8").Value = 14666
$ws.Range("F9").Value = 1378
$ws.Range("F11").Value = 5972
$ws.Range("F12").Value = 1000
$ws.Range("F18").Value = 39
$ws.Range("F19").Value = 841
$ws.Range("F20").Value = 3427
$ws.Range("F22").Value = 2497
$ws.Range("F23").Value = 621
$ws.Range("F31").Value = 1097
$ws.Range("F32").Value = 1473
$ws.Range("F34").Value = 19
$ws.Range("F35").Value = 5057
$ws.Range("F36").Value = 4982
$ws.Range("F38").Value = 694
$ws.Range("F39").Value = 3332
$ws.Range("F42").Value = 130
$ws.Range("F45").Value = 641
$ws.Range("F46").Value = 309
